# Chronologie2023.xlsx - "Add files via upload"
#
# The sheet "Feuil1" is a line-by-line log of anuran (frog/toad) chorus
# observations. This change appends one new observation as row 42:
#   05/05/2023, RASY, N/A, Baie-Saint-Paul, Capitale-Nationale, C/D, Cote 3,
#   "Présence de masses d'oeufs, donnée soumise à l'AARQ", Suzanne Couture

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 42

# --- Formatting -------------------------------------------------------
# Reuse the existing per-column look-and-feel instead of inventing new
# styles: columns A, C, D, E, G, H, I keep the same formatting as the row
# directly above (row 41), while the two "coded" columns (B = species
# code, F = climate-zone code) pick up the colour coding already used
# elsewhere in the sheet for the same codes ("RASY" / "C/D").
$ws.Range("A41").Copy()
$ws.Range("A$newRow").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B4").Copy()                       # an existing "RASY" row
$ws.Range("B$newRow").PasteSpecial(-4122)

$ws.Range("C41").Copy()
$ws.Range("C$newRow").PasteSpecial(-4122)

$ws.Range("D41").Copy()
$ws.Range("D$newRow").PasteSpecial(-4122)

$ws.Range("E41").Copy()
$ws.Range("E$newRow").PasteSpecial(-4122)

$ws.Range("F13").Copy()                      # an existing "C/D" row
$ws.Range("F$newRow").PasteSpecial(-4122)

$ws.Range("G41").Copy()
$ws.Range("G$newRow").PasteSpecial(-4122)

$ws.Range("H41").Copy()
$ws.Range("H$newRow").PasteSpecial(-4122)

$ws.Range("I41").Copy()
$ws.Range("I$newRow").PasteSpecial(-4122)

# --- Values -------------------------------------------------------------
$ws.Range("A$newRow").Value = 45051   # 2023-05-05
$ws.Range("B$newRow").Value = "RASY"
$ws.Range("C$newRow").Value = "N/A"
$ws.Range("D$newRow").Value = "Baie-Saint-Paul"
$ws.Range("E$newRow").Value = "Capitale-Nationale"
$ws.Range("F$newRow").Value = "C/D"
$ws.Range("G$newRow").Value = "Cote 3"
$ws.Range("H$newRow").Value = "Présence de masses d'oeufs, donnée soumise à l'AARQ"
$ws.Range("I$newRow").Value = "Suzanne Couture"

# Move the active selection below the newly-added row, like a user who
# just finished typing the entry would leave it.
$null = $ws.Range("A45").Select()
